{"js": "// Wrap a <w:document>...</w:document> body fragment into the FlatOPC\n// package format required by Range.insertOoxml().\nfunction wrapFlatOpc(bodyXml) {\n  return `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n<pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\" pkg:padding=\"512\">\n<pkg:xmlData>\n<Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\"><Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/></Relationships>\n</pkg:xmlData>\n</pkg:part>\n<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n<pkg:xmlData>\n<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">${bodyXml}</w:document>\n</pkg:xmlData>\n</pkg:part>\n</pkg:package>`;\n}\n\n// Each entry: the paragraph's current plain-text (used to locate it) and\n// the replacement run/proofErr markup (a spell-check pass splitting the\n// text into runs around words Word's proofer flags, plus one appended\n// sentence on the last paragraph).\nconst edits = [\n  {\n    match: \"Installation of jsonwebtoken : npm install jsonwebtoken\",\n    xml:\n      '<w:body><w:p>' +\n      '<w:r><w:t xml:space=\"preserve\">Installation of jsonwebtoken : </w:t></w:r>' +\n      '<w:proofErr w:type=\"spellStart\"/><w:r><w:t>npm</w:t></w:r><w:proofErr w:type=\"spellEnd\"/>' +\n      '<w:r><w:t xml:space=\"preserve\"> install </w:t></w:r>' +\n      '<w:proofErr w:type=\"spellStart\"/><w:r><w:t>jsonwebtoken</w:t></w:r><w:proofErr w:type=\"spellEnd\"/>' +\n      '</w:p></w:body>'\n  },\n  {\n    match:\n      \"Make a file named user.js and copy the code from the study materials.\" +\n      \"hashing is the one way where data can be encrypted but it\\u2019s a one way . datas that has been encrypted cannot be decrypted back.\",\n    xml:\n      '<w:body><w:p>' +\n      '<w:r><w:lastRenderedPageBreak/><w:t xml:space=\"preserve\">Make a file named user.js and copy the code from the study </w:t></w:r>' +\n      '<w:proofErr w:type=\"spellStart\"/><w:r><w:t>materials.</w:t></w:r><w:r><w:t>hashing</w:t></w:r><w:proofErr w:type=\"spellEnd\"/>' +\n      '<w:r><w:t xml:space=\"preserve\"> is the one way where data can be encrypted but it\\u2019s a one way . </w:t></w:r>' +\n      '<w:proofErr w:type=\"spellStart\"/><w:r><w:t>datas</w:t></w:r><w:proofErr w:type=\"spellEnd\"/>' +\n      '<w:r><w:t xml:space=\"preserve\"> that has been encrypted cannot be decrypted back.</w:t></w:r>' +\n      '</w:p></w:body>'\n  },\n  {\n    match: \"In the blog.js set the reference of the user for setting the scema .\",\n    xml:\n      '<w:body><w:p>' +\n      '<w:r><w:t xml:space=\"preserve\">In the blog.js set the reference of the user for setting the </w:t></w:r>' +\n      '<w:proofErr w:type=\"spellStart\"/><w:r><w:t>scema</w:t></w:r><w:proofErr w:type=\"spellEnd\"/>' +\n      '<w:r><w:t xml:space=\"preserve\"> .</w:t></w:r>' +\n      '</w:p></w:body>'\n  },\n  {\n    match:\n      \"Then install the bycrypt to create one way hash. Npm install bycrypt is the command . make a file name users.js in the controller and import it in the app.js file.\",\n    xml:\n      '<w:body><w:p>' +\n      '<w:r><w:t xml:space=\"preserve\">Then install the </w:t></w:r>' +\n      '<w:proofErr w:type=\"spellStart\"/><w:r><w:t>bycrypt</w:t></w:r><w:proofErr w:type=\"spellEnd\"/>' +\n      '<w:r><w:t xml:space=\"preserve\"> to create one way hash. </w:t></w:r>' +\n      '<w:proofErr w:type=\"spellStart\"/><w:r><w:t>Npm</w:t></w:r><w:proofErr w:type=\"spellEnd\"/>' +\n      '<w:r><w:t xml:space=\"preserve\"> install </w:t></w:r>' +\n      '<w:proofErr w:type=\"spellStart\"/><w:r><w:t>bycrypt</w:t></w:r><w:proofErr w:type=\"spellEnd\"/>' +\n      '<w:r><w:t xml:space=\"preserve\"> is the command . make a file name users.js in the controller and import it in the app.js file.</w:t></w:r>' +\n      '<w:r><w:t xml:space=\"preserve\"> Copy the code of users.js from the study materials and export it .</w:t></w:r>' +\n      '</w:p></w:body>'\n  }\n];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nfor (const edit of edits) {\n  let target = null;\n  for (let i = 0; i < paragraphs.items.length; i++) {\n    if (paragraphs.items[i].text === edit.match) {\n      target = paragraphs.items[i];\n      break;\n    }\n  }\n  if (!target) {\n    throw new Error(\"Could not locate paragraph with text: \" + edit.match);\n  }\n  const range = target.getRange();\n  range.insertOoxml(wrapFlatOpc(edit.xml), \"Replace\");\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\nfunction Wrap-FlatOpc([string]$bodyXml) {\n    return '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' + `\n        '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' + `\n        '<pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\" pkg:padding=\"512\">' + `\n        '<pkg:xmlData>' + `\n        '<Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\"><Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/></Relationships>' + `\n        '</pkg:xmlData>' + `\n        '</pkg:part>' + `\n        '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' + `\n        '<pkg:xmlData>' + `\n        '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' + $bodyXml + '</w:document>' + `\n        '</pkg:xmlData>' + `\n        '</pkg:part>' + `\n        '</pkg:package>'\n}\n\n# Each entry: the paragraph's current plain-text (used to locate it) and the\n# replacement run/proofErr markup -- a spell-check pass that splits the text\n# into runs around words Word's proofer flags, plus one appended sentence on\n# the last paragraph.\n$edits = @(\n    @{\n        Match = \"Installation of jsonwebtoken : npm install jsonwebtoken\"\n        Xml = '<w:body><w:p><w:r><w:t xml:space=\"preserve\">Installation of jsonwebtoken : </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>npm</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> install </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>jsonwebtoken</w:t></w:r><w:proofErr w:type=\"spellEnd\"/></w:p></w:body>'\n    },\n    @{\n        Match = \"Make a file named user.js and copy the code from the study materials.hashing is the one way where data can be encrypted but it\" + [char]0x2019 + \"s a one way . datas that has been encrypted cannot be decrypted back.\"\n        Xml = '<w:body><w:p><w:r><w:lastRenderedPageBreak/><w:t xml:space=\"preserve\">Make a file named user.js and copy the code from the study </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>materials.</w:t></w:r><w:r><w:t>hashing</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> is the one way where data can be encrypted but it' + [char]0x2019 + 's a one way . </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>datas</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> that has been encrypted cannot be decrypted back.</w:t></w:r></w:p></w:body>'\n    },\n    @{\n        Match = \"In the blog.js set the reference of the user for setting the scema .\"\n        Xml = '<w:body><w:p><w:r><w:t xml:space=\"preserve\">In the blog.js set the reference of the user for setting the </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>scema</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> .</w:t></w:r></w:p></w:body>'\n    },\n    @{\n        Match = \"Then install the bycrypt to create one way hash. Npm install bycrypt is the command . make a file name users.js in the controller and import it in the app.js file.\"\n        Xml = '<w:body><w:p><w:r><w:t xml:space=\"preserve\">Then install the </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>bycrypt</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> to create one way hash. </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>Npm</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> install </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>bycrypt</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> is the command . make a file name users.js in the controller and import it in the app.js file.</w:t></w:r><w:r><w:t xml:space=\"preserve\"> Copy the code of users.js from the study materials and export it .</w:t></w:r></w:p></w:body>'\n    }\n)\n\nforeach ($edit in $edits) {\n    $paras = $d.Paragraphs\n    $count = $paras.Count\n    for ($i = 1; $i -le $count; $i++) {\n        $p = $paras.Item($i)\n        $t = $p.Range.Text\n        $t = $t.TrimEnd([char]13, [char]7)\n        if ($t -eq $edit.Match) {\n            $rng = $p.Range\n            $rng.MoveEnd(1, -1) | Out-Null\n            $rng.Text = \"\"\n            $rng.InsertXML((Wrap-FlatOpc $edit.Xml))\n            break\n        }\n    }\n}\n"}
